$p = $ppt.ActivePresentation
try {
    $css = $p.ColorSchemes
    Write-Host "ColorSchemes: $css Count=$($css.Count)"
    for ($i=1; $i -le $css.Count; $i++) {
        $cs = $css.Item($i)
        Write-Host "  Item $i : $cs"
    }
} catch {
    Write-Host "ERR: $_"
}
